$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T45")

$ws.Range("B2").Value = -0.0958905286338048
$ws.Range("C2").Value = 1.890977232749403
$ws.Range("D2").Value = 16.63270717711324
$ws.Range("E2").Value = 4.07832161276097
$ws.Range("F2").Value = 4.173141263113404
$ws.Range("G2").Value = 22
$ws.Range("B3").Value = 0.1616684523960579
$ws.Range("C3").Value = 1.751782844777193
$ws.Range("D3").Value = 11.43642227753171
$ws.Range("E3").Value = 3.381777975789024
$ws.Range("F3").Value = 3.461329205447403
$ws.Range("G3").Value = 21
$ws.Range("B4").Value = -0.5325954658642663
$ws.Range("C4").Value = 1.006200763201028
$ws.Range("D4").Value = 4.239438742259294
$ws.Range("E4").Value = 2.058989738259833
$ws.Range("F4").Value = 2.040583201466039
$ws.Range("G4").Value = 20
$ws.Range("B5").Value = 0.08803471659571853
$ws.Range("C5").Value = 0.7012671322584916
$ws.Range("D5").Value = 1.718504603501026
$ws.Range("E5").Value = 1.310917466319305
$ws.Range("F5").Value = 1.343799244086288
$ws.Range("G5").Value = 19
$ws.Range("B6").Value = 0.03921355053569192
$ws.Range("C6").Value = 0.7111885899396917
$ws.Range("D6").Value = 1.651775566409927
$ws.Range("E6").Value = 1.28521421032057
$ws.Range("F6").Value = 1.32185879706788
$ws.Range("G6").Value = 18
$ws.Range("B7").Value = 0.001844256264950899
$ws.Range("C7").Value = 0.5263250456906655
$ws.Range("D7").Value = 0.5109976273144233
$ws.Range("E7").Value = 0.7148409804386031
$ws.Range("F7").Value = 0.7368387646970881
$ws.Range("G7").Value = 17
$ws.Range("B8").Value = 0.08035112964095253
$ws.Range("C8").Value = 0.5395919000237238
$ws.Range("D8").Value = 0.5895437054918792
$ws.Range("E8").Value = 0.7678174949113098
$ws.Range("F8").Value = 0.7886443399199595
$ws.Range("G8").Value = 16
$ws.Range("B9").Value = 0.2273427043598581
$ws.Range("C9").Value = 0.4403849132099676
$ws.Range("D9").Value = 0.3519323016671984
$ws.Range("E9").Value = 0.5932388234658942
$ws.Range("F9").Value = 0.5671806179077581
$ws.Range("G9").Value = 15
$ws.Range("B10").Value = 0.1748763756256576
$ws.Range("C10").Value = 0.409076363660635
$ws.Range("D10").Value = 0.3114511502836552
$ws.Range("E10").Value = 0.5580780861883534
$ws.Range("F10").Value = 0.5499770379433091
$ws.Range("G10").Value = 14
$ws.Range("B11").Value = 0.2035941852220243
$ws.Range("C11").Value = 0.3624616766655595
$ws.Range("D11").Value = 0.1868369991019794
$ws.Range("E11").Value = 0.432246456436579
$ws.Range("F11").Value = 0.396865141766788
